# Updates cryptocurrency price (D) and 1h volume change (E) values
# to match the latest scraped data, preserving text formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "43.780.91"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +2.12%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.334.11"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +1.83%  "
$ws.Range("E4").Value = "  -0.76%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "312.95"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.54%  "
$ws.Range("E6").Value = "  +4.56%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.632"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +1.04%  "
$ws.Range("E8").Value = "  -0.27%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.619"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +2.75%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "41.21"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +4.90%  "
$ws.Range("E11").Value = "  +1.45%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "8.56"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +2.69%  "
$ws.Range("E13").Value = "  -1.00%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "1.01"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +1.81%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "15.51"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +1.97%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "2.685.78"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +1.74%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "2.325.53"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +1.62%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "43.685.82"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +2.11%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "7.55"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +1.64%  "
$ws.Range("E20").Value = "  +1.46%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "13.03"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -4.75%  "
$ws.Range("E22").Value = "  +0.84%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "3.47"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -2.77%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "268.00"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +0.93%  "
$ws.Range("E25").Value = "  +3.53%  "
$ws.Range("E26").Value = "  -0.23%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "7.68"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +9.41%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "11.15"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +3.38%  "
$ws.Range("E29").Value = "  -1.46%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "39.76"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +8.37%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "22.58"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +0.11%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "167.80"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +0.35%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.0897"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +3.25%  "
$ws.Range("E34").Value = "  +8.70%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.132"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +0.88%  "
$ws.Range("E36").Value = "  +2.60%  "
$ws.Range("E37").Value = "  +3.88%  "
$ws.Range("E38").Value = "  +4.02%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "2.92"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +9.84%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "3.80"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +3.20%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "1.72"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +8.66%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "104.63"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +11.26%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "13.74"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +14.36%  "
$ws.Range("E44").Value = "  +4.71%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "71.75"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +2.28%  "
$ws.Range("E46").Value = "  -0.07%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "114.75"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +2.44%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.221"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +18.17%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "1.659.41"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -4.44%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "8.98"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +3.09%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "76.22"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -4.85%  "
